$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-16 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-17 Sunday", 2) | Out-Null
$d.Content.Find.Execute("71-43=28", $true, $false, $false, $false, $false, $true, 1, $false, "82-6=76", 2) | Out-Null
$d.Content.Find.Execute("81-19=62", $true, $false, $false, $false, $false, $true, 1, $false, "56+8=64", 2) | Out-Null
$d.Content.Find.Execute("16+79=95", $true, $false, $false, $false, $false, $true, 1, $false, "29+44=73", 2) | Out-Null
$d.Content.Find.Execute("19+38=57", $true, $false, $false, $false, $false, $true, 1, $false, "67+8=75", 2) | Out-Null
$d.Content.Find.Execute("90-29=61", $true, $false, $false, $false, $false, $true, 1, $false, "75+16=91", 2) | Out-Null
$d.Content.Find.Execute("74-46=28", $true, $false, $false, $false, $false, $true, 1, $false, "56+18=74", 2) | Out-Null
$d.Content.Find.Execute("93-89=4", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("79+13=92", $true, $false, $false, $false, $false, $true, 1, $false, "75+7=82", 2) | Out-Null
$d.Content.Find.Execute("32-26=6", $true, $false, $false, $false, $false, $true, 1, $false, "29+54=83", 2) | Out-Null
$d.Content.Find.Execute("91-14=77", $true, $false, $false, $false, $false, $true, 1, $false, "34-19=15", 2) | Out-Null
$d.Content.Find.Execute("28+9=37", $true, $false, $false, $false, $false, $true, 1, $false, "17+19=36", 2) | Out-Null
$d.Content.Find.Execute("59+25=84", $true, $false, $false, $false, $false, $true, 1, $false, "93-67=26", 2) | Out-Null
$d.Content.Find.Execute("39+53=92", $true, $false, $false, $false, $false, $true, 1, $false, "11-8=3", 2) | Out-Null
$d.Content.Find.Execute("46+17=63", $true, $false, $false, $false, $false, $true, 1, $false, "59+32=91", 2) | Out-Null
$d.Content.Find.Execute("53-29=24", $true, $false, $false, $false, $false, $true, 1, $false, "52-8=44", 2) | Out-Null
$d.Content.Find.Execute("67-19=48", $true, $false, $false, $false, $false, $true, 1, $false, "81-26=55", 2) | Out-Null
$d.Content.Find.Execute("65+18=83", $true, $false, $false, $false, $false, $true, 1, $false, "72-33=39", 2) | Out-Null
$d.Content.Find.Execute("25-16=9", $true, $false, $false, $false, $false, $true, 1, $false, "93-55=38", 2) | Out-Null
$d.Content.Find.Execute("88-79=9", $true, $false, $false, $false, $false, $true, 1, $false, "93-7=86", 2) | Out-Null
$d.Content.Find.Execute("13-7=6", $true, $false, $false, $false, $false, $true, 1, $false, "67+6=73", 2) | Out-Null
$d.Content.Find.Execute("81-72=9", $true, $false, $false, $false, $false, $true, 1, $false, "40-28=12", 2) | Out-Null
$d.Content.Find.Execute("40-32=8", $true, $false, $false, $false, $false, $true, 1, $false, "19+32=51", 2) | Out-Null
$d.Content.Find.Execute("80-17=63", $true, $false, $false, $false, $false, $true, 1, $false, "75-48=27", 2) | Out-Null
$d.Content.Find.Execute("16+29=45", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=5", 2) | Out-Null
$d.Content.Find.Execute("9+39=48", $true, $false, $false, $false, $false, $true, 1, $false, "57+26=83", 2) | Out-Null
$d.Content.Find.Execute("19+63=82", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=59", 2) | Out-Null
$d.Content.Find.Execute("8+56=64", $true, $false, $false, $false, $false, $true, 1, $false, "19+23=42", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "95-16=79", 2) | Out-Null
$d.Content.Find.Execute("8+69=77", $true, $false, $false, $false, $false, $true, 1, $false, "7+79=86", 2) | Out-Null
$d.Content.Find.Execute("79+16=95", $true, $false, $false, $false, $false, $true, 1, $false, "15+59=74", 2) | Out-Null
$d.Content.Find.Execute("38+15=53", $true, $false, $false, $false, $false, $true, 1, $false, "46+36=82", 2) | Out-Null
$d.Content.Find.Execute("29+39=68", $true, $false, $false, $false, $false, $true, 1, $false, "28+23=51", 2) | Out-Null
$d.Content.Find.Execute("58-49=9", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=86", 2) | Out-Null
$d.Content.Find.Execute("39+54=93", $true, $false, $false, $false, $false, $true, 1, $false, "58+38=96", 2) | Out-Null
$d.Content.Find.Execute("45-36=9", $true, $false, $false, $false, $false, $true, 1, $false, "92-89=3", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "38+58=96", 2) | Out-Null
$d.Content.Find.Execute("71-69=2", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=28", 2) | Out-Null
$d.Content.Find.Execute("63-16=47", $true, $false, $false, $false, $false, $true, 1, $false, "44-6=38", 2) | Out-Null
$d.Content.Find.Execute("35+7=42", $true, $false, $false, $false, $false, $true, 1, $false, "16+5=21", 2) | Out-Null
$d.Content.Find.Execute("95-6=89", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("94-7=87", $true, $false, $false, $false, $false, $true, 1, $false, "62-44=18", 2) | Out-Null
$d.Content.Find.Execute("81-57=24", $true, $false, $false, $false, $false, $true, 1, $false, "66-37=29", 2) | Out-Null
$d.Content.Find.Execute("81-66=15", $true, $false, $false, $false, $false, $true, 1, $false, "47+49=96", 2) | Out-Null
$d.Content.Find.Execute("64+19=83", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=59", 2) | Out-Null
$d.Content.Find.Execute("82-3=79", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=52", 2) | Out-Null
$d.Content.Find.Execute("61-37=24", $true, $false, $false, $false, $false, $true, 1, $false, "19+46=65", 2) | Out-Null
$d.Content.Find.Execute("3+69=72", $true, $false, $false, $false, $false, $true, 1, $false, "31-18=13", 2) | Out-Null
$d.Content.Find.Execute("81-14=67", $true, $false, $false, $false, $false, $true, 1, $false, "29+68=97", 2) | Out-Null
$d.Content.Find.Execute("17+54=71", $true, $false, $false, $false, $false, $true, 1, $false, "38+38=76", 2) | Out-Null
$d.Content.Find.Execute("74-25=49", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=55", 2) | Out-Null
$d.Content.Find.Execute("16+56=72", $true, $false, $false, $false, $false, $true, 1, $false, "27+67=94", 2) | Out-Null
$d.Content.Find.Execute("84-67=17", $true, $false, $false, $false, $false, $true, 1, $false, "40-6=34", 2) | Out-Null
$d.Content.Find.Execute("4+18=22", $true, $false, $false, $false, $false, $true, 1, $false, "8+83=91", 2) | Out-Null
$d.Content.Find.Execute("40-27=13", $true, $false, $false, $false, $false, $true, 1, $false, "5+79=84", 2) | Out-Null
$d.Content.Find.Execute("61-53=8", $true, $false, $false, $false, $false, $true, 1, $false, "6+27=33", 2) | Out-Null
$d.Content.Find.Execute("36-29=7", $true, $false, $false, $false, $false, $true, 1, $false, "82-77=5", 2) | Out-Null
$d.Content.Find.Execute("59+3=62", $true, $false, $false, $false, $false, $true, 1, $false, "7+69=76", 2) | Out-Null
$d.Content.Find.Execute("18+45=63", $true, $false, $false, $false, $false, $true, 1, $false, "48+27=75", 2) | Out-Null
$d.Content.Find.Execute("68+13=81", $true, $false, $false, $false, $false, $true, 1, $false, "60-34=26", 2) | Out-Null
$d.Content.Find.Execute("92-55=37", $true, $false, $false, $false, $false, $true, 1, $false, "14+58=72", 2) | Out-Null
$d.Content.Find.Execute("72-57=15", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("74+8=82", $true, $false, $false, $false, $false, $true, 1, $false, "45+48=93", 2) | Out-Null
$d.Content.Find.Execute("40-26=14", $true, $false, $false, $false, $false, $true, 1, $false, "81-58=23", 2) | Out-Null
$d.Content.Find.Execute("42+49=91", $true, $false, $false, $false, $false, $true, 1, $false, "40-7=33", 2) | Out-Null
$d.Content.Find.Execute("27+48=75", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=6", 2) | Out-Null
$d.Content.Find.Execute("54-29=25", $true, $false, $false, $false, $false, $true, 1, $false, "30-8=22", 2) | Out-Null
$d.Content.Find.Execute("33-29=4", $true, $false, $false, $false, $false, $true, 1, $false, "15+48=63", 2) | Out-Null
$d.Content.Find.Execute("74-17=57", $true, $false, $false, $false, $false, $true, 1, $false, "46+46=92", 2) | Out-Null
$d.Content.Find.Execute("13+48=61", $true, $false, $false, $false, $false, $true, 1, $false, "15+26=41", 2) | Out-Null
$d.Content.Find.Execute("76-9=67", $true, $false, $false, $false, $false, $true, 1, $false, "46+27=73", 2) | Out-Null
$d.Content.Find.Execute("80-62=18", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=22", 2) | Out-Null
$d.Content.Find.Execute("47+9=56", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=72", 2) | Out-Null
$d.Content.Find.Execute("60-4=56", $true, $false, $false, $false, $false, $true, 1, $false, "18+6=24", 2) | Out-Null
$d.Content.Find.Execute("93-17=76", $true, $false, $false, $false, $false, $true, 1, $false, "36-7=29", 2) | Out-Null
$d.Content.Find.Execute("84-17=67", $true, $false, $false, $false, $false, $true, 1, $false, "28+17=45", 2) | Out-Null
$d.Content.Find.Execute("92-45=47", $true, $false, $false, $false, $false, $true, 1, $false, "33+59=92", 2) | Out-Null
$d.Content.Find.Execute("84-65=19", $true, $false, $false, $false, $false, $true, 1, $false, "93-16=77", 2) | Out-Null
$d.Content.Find.Execute("80-57=23", $true, $false, $false, $false, $false, $true, 1, $false, "91-65=26", 2) | Out-Null
$d.Content.Find.Execute("4+27=31", $true, $false, $false, $false, $false, $true, 1, $false, "14-5=9", 2) | Out-Null
$d.Content.Find.Execute("53-8=45", $true, $false, $false, $false, $false, $true, 1, $false, "58+9=67", 2) | Out-Null
$d.Content.Find.Execute("90-33=57", $true, $false, $false, $false, $false, $true, 1, $false, "63-44=19", 2) | Out-Null
$d.Content.Find.Execute("55+9=64", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2) | Out-Null
$d.Content.Find.Execute("77-8=69", $true, $false, $false, $false, $false, $true, 1, $false, "4+47=51", 2) | Out-Null
$d.Content.Find.Execute("54-17=37", $true, $false, $false, $false, $false, $true, 1, $false, "25-18=7", 2) | Out-Null
$d.Content.Find.Execute("82-65=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+47=53", 2) | Out-Null
$d.Content.Find.Execute("78+13=91", $true, $false, $false, $false, $false, $true, 1, $false, "18+67=85", 2) | Out-Null
$d.Content.Find.Execute("82-29=53", $true, $false, $false, $false, $false, $true, 1, $false, "35+8=43", 2) | Out-Null
$d.Content.Find.Execute("90-72=18", $true, $false, $false, $false, $false, $true, 1, $false, "94-8=86", 2) | Out-Null
$d.Content.Find.Execute("89+9=98", $true, $false, $false, $false, $false, $true, 1, $false, "19+66=85", 2) | Out-Null
$d.Content.Find.Execute("7+86=93", $true, $false, $false, $false, $false, $true, 1, $false, "49+44=93", 2) | Out-Null
$d.Content.Find.Execute("18+43=61", $true, $false, $false, $false, $false, $true, 1, $false, "91-6=85", 2) | Out-Null
$d.Content.Find.Execute("65+28=93", $true, $false, $false, $false, $false, $true, 1, $false, "74+7=81", 2) | Out-Null
$d.Content.Find.Execute("8+58=66", $true, $false, $false, $false, $false, $true, 1, $false, "2+19=21", 2) | Out-Null
$d.Content.Find.Execute("62-56=6", $true, $false, $false, $false, $false, $true, 1, $false, "86-79=7", 2) | Out-Null
$d.Content.Find.Execute("25+6=31", $true, $false, $false, $false, $false, $true, 1, $false, "14+58=72", 2) | Out-Null
$d.Content.Find.Execute("16-9=7", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("6+56=62", $true, $false, $false, $false, $false, $true, 1, $false, "25+49=74", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "77+6=83", 2) | Out-Null
$d.Content.Find.Execute("49+4=53", $true, $false, $false, $false, $false, $true, 1, $false, "59+13=72", 2) | Out-Null
$d.Content.Find.Execute("26+65=91", $true, $false, $false, $false, $false, $true, 1, $false, "49+23=72", 2) | Out-Null

Write-Output "Replaced 101 items"
